# Scrub sample/live account data from the REPO and DEPO sheets,
# replacing it with generic placeholder values ("TEST" / all-zero
# account number) before the template is shared publicly.
#
# NOTE: "TEST" is written before the zero-account-number so that the
# new shared-string table entries are created in the same order as in
# the target workbook (TEST first, then the zero account number).
$wb = $excel.ActiveWorkbook

# --- REPO sheet ---
$wsRepo = $wb.Worksheets.Item("REPO")
$wsRepo.Range("B2").Value = "TEST"
# Leading apostrophe keeps this a genuine text entry (column A is
# formatted/quote-prefixed as Text) instead of letting Excel reinterpret
# the all-numeric string and silently create a new number-format style.
$wsRepo.Range("A2").Value = "'000000000000000"
$wsRepo.Range("C2").Value = "TEST"
$wsRepo.Range("G2").Value = "TEST"
$wsRepo.Range("J2").Value = "TEST"
$wsRepo.Range("K2").Value = "TEST"
$wsRepo.Range("A2").Select()

# --- DEPO sheet ---
$wsDepo = $wb.Worksheets.Item("DEPO")
$wsDepo.Range("B2").Value = "TEST"
$wsDepo.Range("A2").Value = "'000000000000000"
$wsDepo.Range("C2").Value = "TEST"
$wsDepo.Activate()
$wsDepo.Range("F18").Select()
